$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert two new rows for the "GetData Service" workblock entry, right before
# the old row 25 (wbProcessTransaction_Type). This shifts everything below
# (including the Task1/2/3 block) down by two rows.
$ws.Rows("25:26").Insert()

# Restore row height / custom-height formatting that Insert() does not carry
# over to the freshly inserted rows.
$ws.Rows("25:26").RowHeight = $ws.Rows("24").RowHeight

# Row 25 should look like the plain "Type" rows (no special style on column B)
$ws.Cells.Item(25, 2).Style = "Normal"

# Populate the new workblock definition
$ws.Cells.Item(25, 1).Value = "wbGetDataTask_Type"
$ws.Cells.Item(25, 2).Value = "Main, Framework, GetData Service"
$ws.Cells.Item(25, 3).Value = "Name of Workblock"
$ws.Cells.Item(26, 1).Value = "wbGetDataTask_SuppressSuccessful"
$ws.Cells.Item(26, 2).Value = $true
$ws.Cells.Item(26, 3).Value = "Do not log successful executions of wb"

# Task1 (FirstRun) is renamed to "FirstRunService" and disabled by default.
$ws.Cells.Item(38, 2).Value = "FirstRunService"
$ws.Cells.Item(40, 2).Value = $false

# Task2 is now assigned to the new "GetDataService" task.
$ws.Cells.Item(41, 2).Value = "GetDataService"

# Update the active selection to match the author's saved view state.
$ws.Range("B41").Select() | Out-Null
